$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("max-arrecad")

$values = @(
    "midia_independente",
    "disputa",
    "herois",
    "jogos",
    "erotismo",
    "terror",
    "politica",
    "religiosidade",
    "humor",
    "nenhuma",
    "folclore",
    "ficcao_cientifica",
    "lgbtqiamais",
    "fiq",
    "questoes_genero",
    "ccxp",
    "angelo_agostini",
    "hqmix",
    "zine",
    "fantasia",
    "webformatos",
    "saloes_humor"
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row++
}
